$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.637.67"
$ws.Range("E2").Value = "  -0.28%  "

# Row 3
$ws.Range("D3").Value = "1.642.30"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.26"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6
$ws.Range("E6").Value = "  +1.30%  "

# Row 7
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +0.77%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.27"
$ws.Range("E10").Value = "  +0.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("E11").Value = "  -0.05%  "

# Row 12
$ws.Range("D12").Value = "1.870.13"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.645.29"
$ws.Range("E13").Value = "  +0.20%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.20"
$ws.Range("E14").Value = "  +2.31%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.530"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.39"
$ws.Range("E16").Value = "  +2.84%  "

# Row 17
$ws.Range("D17").Value = "26.681.93"
$ws.Range("E17").Value = "  +0.02%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.64"
$ws.Range("E19").Value = "  -0.82%  "

# Row 20
$ws.Range("E20").Value = "  +0.17%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  +1.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.31"
$ws.Range("E22").Value = "  +2.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.52"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +12.65%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.51"
$ws.Range("E25").Value = "  -1.52%  "

# Row 26
$ws.Range("E26").Value = "  +0.30%  "

# Row 27
$ws.Range("E27").Value = "  -1.01%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("E28").Value = "  +4.45%  "

# Row 29
$ws.Range("E29").Value = "  +1.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0517"
$ws.Range("E30").Value = "  +2.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  +0.34%  "

# Row 32
$ws.Range("E32").Value = "  +2.20%  "

# Row 33
$ws.Range("E33").Value = "  +1.82%  "

# Row 34
$ws.Range("D34").Value = "1.277.05"
$ws.Range("E34").Value = "  +4.21%  "

# Row 35
$ws.Range("E35").Value = "  +2.73%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0181"
$ws.Range("E36").Value = "  +5.13%  "

# Row 37
$ws.Range("E37").Value = "  +0.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("E38").Value = "  +6.60%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("E39").Value = "  +2.60%  "

# Row 40
$ws.Range("E40").Value = "  +0.18%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.818"
$ws.Range("E41").Value = "  +2.73%  "

# Row 42
$ws.Range("E42").Value = "  -1.46%  "

# Row 43
$ws.Range("E43").Value = "  +2.30%  "

# Row 44
$ws.Range("D44").Value = "1.781.16"
$ws.Range("E44").Value = "  +0.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.05"
$ws.Range("E45").Value = "  -0.74%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.89"
$ws.Range("E46").Value = "  +8.24%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("E47").Value = "  +2.05%  "

# Row 48
$ws.Range("E48").Value = "  +0.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("E49").Value = "  +1.97%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  +2.88%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.406"
$ws.Range("E51").Value = "  -0.53%  "
